$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.751.88'
$ws.Range("E2").Value = '  -2.86%  '

$ws.Range("D3").Value = '1.780.88'
$ws.Range("E3").Value = '  -3.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.81'
$ws.Range("E5").Value = '  -6.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5029'
$ws.Range("E7").Value = '  -4.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.36'
$ws.Range("E8").Value = '  -5.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2498'
$ws.Range("E9").Value = '  -21.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06200'
$ws.Range("E10").Value = '  -8.67%  '

$ws.Range("D11").Value = '1.803.49'
$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06772'
$ws.Range("E12").Value = '  -12.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.79'
$ws.Range("E13").Value = '  -21.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6139'
$ws.Range("E14").Value = '  -21.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '79.05'
$ws.Range("E15").Value = '  -10.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.384'
$ws.Range("E16").Value = '  -12.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.003'
$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("D19").Value = '25.765.75'
$ws.Range("E19").Value = '  -2.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.35'
$ws.Range("E20").Value = '  -17.95%  '

$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006374'
$ws.Range("E21").Value = '  -19.84%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.029.31'
$ws.Range("E22").Value = '  -2.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.922'
$ws.Range("E23").Value = '  -15.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.205'
$ws.Range("E24").Value = '  -12.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.049'
$ws.Range("E25").Value = '  -13.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '131.55'
$ws.Range("E26").Value = '  -7.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.904'
$ws.Range("E27").Value = '  -13.51%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.392'
$ws.Range("E28").Value = '  -17.12%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.51'
$ws.Range("E29").Value = '  -14.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '98.92'
$ws.Range("E30").Value = '  -11.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08237'
$ws.Range("E31").Value = '  -5.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.587'
$ws.Range("E32").Value = '  -13.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04314'
$ws.Range("E33").Value = '  -11.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.727'
$ws.Range("E34").Value = '  -4.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.160'
$ws.Range("E35").Value = '  -22.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.036'
$ws.Range("E36").Value = '  -8.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6232'
$ws.Range("E37").Value = '  -14.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.777'
$ws.Range("E38").Value = '  -10.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.117'
$ws.Range("E39").Value = '  -5.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.004'
$ws.Range("E40").Value = '  +0.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.41'
$ws.Range("E41").Value = '  -7.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01463'
$ws.Range("E42").Value = '  -16.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7838'
$ws.Range("E43").Value = '  -12.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3893'
$ws.Range("E44").Value = '  -19.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.198'
$ws.Range("E45").Value = '  -12.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.197'
$ws.Range("E46").Value = '  -18.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05241'
$ws.Range("E47").Value = '  -10.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '52.65'
$ws.Range("E48").Value = '  -11.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.476'
$ws.Range("E50").Value = '  -16.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '29.21'
$ws.Range("E51").Value = '  -16.27%  '
